# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" positioned right before the "总计"
#    sheet (so tab order becomes 2021-Q3, 2021-Q4, 2022-Q1, 总计), carrying
#    the per-fund holdings table for the new quarter.
# 2. Insert a new summary row at the top of "总计"'s data (after the header)
#    for "2022-Q1", pushing the existing 2021-Q4 / 2021-Q3 rows down by one.

$wb = $excel.ActiveWorkbook

$existing = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# --- 1. New "2022-Q1" sheet, inserted before "总计" -------------------------

$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# Match the page margins used across the rest of this workbook (0.75" /
# 0.75" / 1" / 1" / 0.5" / 0.5", i.e. 54/54/72/72/36/36 points) instead of
# the blank-sheet Excel defaults.
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# NOTE: worksheet references here track *position*, not a fixed sheet
# identity. Inserting a sheet before "总计" shifts its index, so the old
# $totalSheet handle now resolves to the newly inserted sheet instead of
# "总计" — re-fetch it by name before touching it again.
$totalSheet = $wb.Worksheets.Item("总计")

# Clone the header-row look (border/bold/center style) from an existing
# quarter sheet so the new header cells share style index 2, instead of
# minting a fresh style.
$existing.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$existing.Range("A2").Copy()
$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'006923"
$newSheet.Range("C2").Value = "前海开源沪港深非周期性行业股票A"
$newSheet.Range("D2").Value = "'0.54"
$newSheet.Range("E2").Value = "'93.77"
$newSheet.Range("F2").Value = "'5.22"
$newSheet.Range("G2").Value = "'0.0282"
$newSheet.Range("H2").Value = 8

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "'006924"
$newSheet.Range("C3").Value = "前海开源沪港深非周期性行业股票C"
$newSheet.Range("D3").Value = "'0.22"
$newSheet.Range("E3").Value = "'93.77"
$newSheet.Range("F3").Value = "'5.22"
$newSheet.Range("G3").Value = "'0.0115"
$newSheet.Range("H3").Value = 8

# The fund-code / numeric-looking text cells above were entered with a
# leading apostrophe so Excel stores them as text rather than numbers;
# strip the resulting quote-prefix formatting so the cells end up with no
# explicit style, matching the rest of the data rows.
$newSheet.Range("B2:B3").ClearFormats()
$newSheet.Range("D2:G3").ClearFormats()

# --- 2. Add a "2022-Q1" summary row to "总计", shifting old rows down ------

$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 2
$totalSheet.Range("D2").Value = 0.04

$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 14
$totalSheet.Range("D3").Value = 11.35

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 14
$totalSheet.Range("D4").Value = 9.640000000000001

# Restore the originally active tab ("2021-Q3") — none of the edits above
# are meant to change which sheet is in focus when the workbook is opened.
$wb.Worksheets.Item("2021-Q3").Select()
